# Daten aktualisiert am 2024-03-17
# - Normalise "Sector" column text to sentence case (only first word/acronym
#   capitalised, e.g. "Financial Services" -> "Financial services").
# - Replace "Endeavour Mining" (EDV) with "EasyJet" (EZJ); "Entain" (ENT)
#   moves up one row and its sector is corrected from "Mining" to
#   "Travel & leisure".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('C2').Value = 'Financial services'
$ws.Range('C4').Value = 'Telecommunications services'
$ws.Range('C7').Value = 'Support services'
$ws.Range('C8').Value = 'Food & tobacco'
$ws.Range('C9').Value = 'Pharmaceuticals & biotechnology'
$ws.Range('C11').Value = 'Life insurance'
$ws.Range('C13').Value = 'Aerospace & defence'
$ws.Range('C15').Value = 'Household goods & home construction'
$ws.Range('C17').Value = 'Household goods & home construction'
$ws.Range('C18').Value = 'Oil & gas producers'
$ws.Range('C20').Value = 'Telecommunications services'
$ws.Range('C21').Value = 'Support services'
$ws.Range('C22').Value = 'Personal goods'
$ws.Range('C23').Value = 'Multiline utilities'
$ws.Range('C25').Value = 'Support services'
$ws.Range('C26').Value = 'Health care equipment & supplies'
$ws.Range('C28').Value = 'Support services'
$ws.Range('C30').Value = 'Industrial Support services'

# EDV/Endeavour Mining removed; ENT/Entain shifts up from row 32 to 31 (and
# its sector is corrected); EZJ/EasyJet is newly inserted at row 32.
$ws.Range('A31').Value = 'ENT'
$ws.Range('B31').Value = 'Entain'
$ws.Range('C31').Value = 'Travel & leisure'
$ws.Range('A32').Value = 'EZJ'
$ws.Range('B32').Value = 'EasyJet'
$ws.Range('C32').Value = 'Travel & leisure'

$ws.Range('C33').Value = 'Support services'
$ws.Range('C34').Value = 'Financial services'
$ws.Range('C35').Value = 'Travel & leisure'
$ws.Range('C39').Value = 'Pharmaceuticals & biotechnology'
$ws.Range('C40').Value = 'Pharmaceuticals & biotechnology'
$ws.Range('C41').Value = 'Electronic equipment & parts'
$ws.Range('C42').Value = 'Pharmaceuticals & biotechnology'
$ws.Range('C43').Value = 'Homebuilding & construction supplies'
$ws.Range('C45').Value = 'Travel & leisure'
$ws.Range('C46').Value = 'Machinery, tools, heavy vehicles, trains & ships'
$ws.Range('C49').Value = 'Financial services'
$ws.Range('C50').Value = 'Travel & leisure'
$ws.Range('C51').Value = 'Support services'
$ws.Range('C52').Value = 'General retailers'
$ws.Range('C54').Value = 'Real estate investment trusts'
$ws.Range('C55').Value = 'Life insurance'
$ws.Range('C57').Value = 'Financial services'
$ws.Range('C58').Value = 'Financial services'
$ws.Range('C59').Value = 'Food & drug retailing'
$ws.Range('C60').Value = 'Aerospace & defence'
$ws.Range('C61').Value = 'Containers & packaging'
$ws.Range('C62').Value = 'Multiline utilities'
$ws.Range('C64').Value = 'General retailers'
$ws.Range('C65').Value = 'Food & drug retailers'
$ws.Range('C67').Value = 'Financial services'
$ws.Range('C68').Value = 'Household goods & home construction'
$ws.Range('C69').Value = 'Life insurance'
$ws.Range('C70').Value = 'Life insurance'
$ws.Range('C71').Value = 'Household goods & home construction'
$ws.Range('C73').Value = 'Support services'
$ws.Range('C76').Value = 'Aerospace & defence'
$ws.Range('C78').Value = 'Software & computer services'
$ws.Range('C79').Value = 'Food & drug retailing'
$ws.Range('C80').Value = 'Financial services'
$ws.Range('C81').Value = 'Collective investments'
$ws.Range('C82').Value = 'Real estate investment trusts'
$ws.Range('C83').Value = 'Multiline utilities'
$ws.Range('C84').Value = 'Oil & gas producers'
$ws.Range('C85').Value = 'General industrials'
$ws.Range('C86').Value = 'General industrials'
$ws.Range('C87').Value = 'Health care equipment & supplies'
$ws.Range('C88').Value = 'General industrials'
$ws.Range('C89').Value = 'Industrial engineering'
$ws.Range('C90').Value = 'Electrical utilities & independent power producers'
$ws.Range('C92').Value = 'Financial services'
$ws.Range('C93').Value = 'Household goods & home construction'
$ws.Range('C94').Value = 'Food & drug retailing'
$ws.Range('C95').Value = 'Personal goods'
$ws.Range('C96').Value = 'Multiline utilities'
$ws.Range('C97').Value = 'Real estate investment trusts'
$ws.Range('C98').Value = 'Mobile telecommunications'
$ws.Range('C99').Value = 'Industrial goods and services'
